$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 24.71000000000042
$ws.Range("H2").Value = [double]"1.460819769243627e-16"
$ws.Range("K2").Value = 41.45922963110626
$ws.Range("L2").Value = "[35.87066008132341, 47.047799180889115]"
$ws.Range("O2").Value = 1.402552876377426
$ws.Range("P2").Value = "[1.2641844311742707, 1.5409213215805808]"
$ws.Range("S2").Value = 54.44639596439006
$ws.Range("T2").Value = "[51.024818597407645, 57.86797333137248]"
$ws.Range("W2").Value = 19.19415415415448
$ws.Range("X2").Value = 18.64998998999031
$ws.Range("Y2").Value = 19.73831831831866

# Row 3 updates
$ws.Range("E3").Value = 22.53000000000008
$ws.Range("H3").Value = [double]"1.460819769243627e-16"
$ws.Range("K3").Value = 48.17500749313863
$ws.Range("L3").Value = "[38.047656563342976, 58.30235842293428]"
$ws.Range("O3").Value = -2.352263568453619
$ws.Range("P3").Value = "[-2.5535267614763892, -2.151000375430849]"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 55.31805044670637
$ws.Range("T3").Value = "[50.18450236253092, 60.45159853088181]"
$ws.Range("W3").Value = 8.434654654654686
$ws.Range("X3").Value = 7.712972972973
$ws.Range("Y3").Value = 9.156336336336373
